# Update the threshold table on Sheet1 with the new measured values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 5.6    # alpha_distance_range Min
$ws.Range("C2").Value = 11.7   # alpha_distance_range Max
$ws.Range("B3").Value = 5.7    # beta_distance_range Min
$ws.Range("C3").Value = 10.4   # beta_distance_range Max
$ws.Range("B4").Value = 0.95   # ratio_threshold_range Min
$ws.Range("C5").Value = 17     # pie_threshold_range Max
